# Rename the header row and add a "Year" column (all 2010) to Sheet1,
# matching the "change column names and year" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# A1: rank -> Rank
# B1: city -> City Name
# C1: dollars_billion -> Overnight International Visitor Spend (US$ bn)
# D1: (new) Year
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitor Spend (US$ bn)"
$ws.Range("D1").Value = "Year"

# --- Data rows ----------------------------------------------------------
# Every data row (2-21) gets a Year value of 2010 in column D
# (D3 previously held a stray blank string; it is overwritten here too).
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = 2010
}

# --- Column widths / selection ------------------------------------------
# (A/B newly sized to fit rank/city text, C widened to fit the long header)
$ws.Columns.Item(1).ColumnWidth = 4.0
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(3).ColumnWidth = 35.75

$ws.Range("D2:D21").Select()
